# Remove post row 750 ("「無限の宇宙を旅した光」...") from the posts sheet.
# All subsequent rows shift up by one (Excel renumbers cell references
# automatically), and the sheet's used range shrinks from 828 to 827 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(750).Delete()
